$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row below the existing data row (16) so the signature rows
# (old 21/22) shift down to 22/23, matching the target layout.
$ws.Rows.Item(17).Insert(-4121)  # xlShiftDown

# Copy formatting (fonts, fills, borders, number formats) from row 16 into
# the freshly-inserted row 17 so the new data row looks identical to the
# existing one.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new "2508" period row for the same worker.
$ws.Range("B17").Value2 = $ws.Range("B16").Value2
$ws.Range("C17").Value2 = $ws.Range("C16").Value2
$ws.Range("D17").Value2 = $ws.Range("D16").Value2
$ws.Range("E17").Value2 = "2508"
$ws.Range("F17").Value2 = 56940
$ws.Range("G17").Value2 = $ws.Range("G16").Value2

# Update the summary "Valor Mora" total (sum of both period rows) and the
# "Cant. Periodos" counter now that a second period has been added.
$ws.Range("E11").Value2 = 70226
$ws.Range("F13").Value2 = 2
